$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.467.43"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "1.826.54"
$ws.Range("E3").Value = "  -1.93%  "
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").Value = "332.45"
$ws.Range("E5").Value = "  -0.55%  "
$ws.Range("E6").Value = "  -0.52%  "
$ws.Range("D7").Value = "0.4580"
$ws.Range("E7").Value = "  -2.17%  "
$ws.Range("D8").Value = "0.3806"
$ws.Range("E8").Value = "  -2.49%  "
$ws.Range("D9").Value = "46.42"
$ws.Range("E9").Value = "  +1.88%  "
$ws.Range("D10").Value = "0.07881"
$ws.Range("E10").Value = "  -1.27%  "
$ws.Range("D11").Value = "0.9697"
$ws.Range("E11").Value = "  -3.15%  "
$ws.Range("D12").Value = "21.03"
$ws.Range("E12").Value = "  -3.14%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "5.884"
$ws.Range("E13").Value = "  -1.73%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.826.02"
$ws.Range("E14").Value = "  -2.17%  "
$ws.Range("D15").Value = "7.057"
$ws.Range("E15").Value = "  -2.44%  "
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("D17").Value = "89.89"
$ws.Range("D18").Value = "0.06645"
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("D20").Value = "17.13"
$ws.Range("E20").Value = "  +0.76%  "
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("D22").Value = "27.440.85"
$ws.Range("E22").Value = "  -0.53%  "
$ws.Range("D23").Value = "5.344"
$ws.Range("E23").Value = "  -1.87%  "
$ws.Range("E24").Value = "  -0.54%  "
$ws.Range("D25").Value = "2.315"
$ws.Range("E25").Value = "  +0.42%  "
$ws.Range("D26").Value = "2.033.48"
$ws.Range("E26").Value = "  -3.21%  "
$ws.Range("D27").Value = "155.35"
$ws.Range("E27").Value = "  -2.67%  "
$ws.Range("E28").Value = "  -2.18%  "
$ws.Range("E29").Value = "  -3.63%  "
$ws.Range("D30").Value = "5.275"
$ws.Range("E30").Value = "  -2.35%  "
$ws.Range("D31").Value = "118.42"
$ws.Range("E31").Value = "  -2.42%  "
$ws.Range("D32").Value = "0.9436"
$ws.Range("E32").Value = "  -3.17%  "
$ws.Range("D33").Value = "0.09310"
$ws.Range("E33").Value = "  -1.53%  "
$ws.Range("D34").Value = "3.596"
$ws.Range("E34").Value = "  -0.64%  "
$ws.Range("D35").Value = "5.246"
$ws.Range("E35").Value = "  -0.74%  "
$ws.Range("D36").Value = "1.319"
$ws.Range("E36").Value = "  -0.90%  "
$ws.Range("D37").Value = "0.05932"
$ws.Range("E37").Value = "  -1.88%  "
$ws.Range("D38").Value = "0.02185"
$ws.Range("E38").Value = "  -1.93%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").Value = "8.047"
$ws.Range("E39").Value = "  -2.64%  "
$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "1.155"
$ws.Range("E40").Value = "  -3.19%  "
$ws.Range("D41").Value = "0.5761"
$ws.Range("E41").Value = "  -2.82%  "
$ws.Range("D42").Value = "0.1828"
$ws.Range("E42").Value = "  -2.79%  "
$ws.Range("D43").Value = "9.983"
$ws.Range("E43").Value = "  -2.32%  "
$ws.Range("D44").Value = "1.259"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").Value = "11.97"
$ws.Range("E45").Value = "  -1.08%  "
$ws.Range("D46").Value = "0.5444"
$ws.Range("E46").Value = "  -3.26%  "
$ws.Range("E47").Value = "  -2.88%  "
$ws.Range("D48").Value = "110.85"
$ws.Range("E48").Value = "  -1.45%  "
$ws.Range("D49").Value = "0.06603"
$ws.Range("E49").Value = "  -2.39%  "
$ws.Range("E50").Value = "  -0.65%  "
$ws.Range("D51").Value = "1.042"
$ws.Range("E51").Value = "  -1.16%  "
